# "added on 19sep 2017" - adds a new "Login" worksheet holding a UID/PWD
# credential pair, and restores the prior cell selections on the existing
# sheets (Sheet1's selection moves to A2; Sheet3 becomes the active tab).

$wb = $excel.ActiveWorkbook

# Sheet1 keeps its data untouched, only its remembered selection changes.
$sheet1 = $wb.Worksheets.Item("Sheet1")
[void]$sheet1.Range("A2").Select()

# Add the new "Login" sheet as the last tab in the workbook.
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "Login"
$ws.Range("A1").Value = "UID"
$ws.Range("B1").Value = "PWD"
$ws.Range("A2").Value = "Ajenkins"
$ws.Range("B2").Value = "Acushnet#1"
[void]$ws.Range("B3").Select()

# Sheet3 becomes the active/selected tab in the workbook.
$sheet3 = $wb.Worksheets.Item("Sheet3")
[void]$sheet3.Activate()
